# Updates cryptos list values (Price column D, Volume(1h) column E)
# per the commit diff. D-column values are forced to remain plain
# text (matching the original inlineStr cells) even when they look
# numeric, by temporarily setting NumberFormat to "@" (Text) before
# assigning the value, then restoring the "Normal" style so no
# stray formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.989.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.27%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.315.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +2.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.336.78"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.104"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.154"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.344"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.11%  "
$ws.Range("E14").Value = "  +4.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.751.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.853.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.327.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.88%  "
$ws.Range("E19").Value = "  +3.28%  "
$ws.Range("E20").Value = "  +3.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.995"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  +6.66%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.66%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("E29").Value = "  +12.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0736"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.04%  "
$ws.Range("E31").Value = "  +5.30%  "
$ws.Range("E32").Value = "  +4.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.26%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.948"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  +6.07%  "
$ws.Range("E39").Value = "  +9.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.57"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.78%  "
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "279.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +16.01%  "
$ws.Range("E45").Value = "  +7.64%  "
$ws.Range("E46").Value = "  +3.79%  "
$ws.Range("E47").Value = "  +4.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.562"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("E50").Value = "  +6.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.19%  "

Write-Host "Applied 75 cell updates (cryptos list refresh)"
